$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above row 2 (shifts old rows 2..22 down to 3..23)
$ws.Rows.Item(2).Insert()

# 2) Populate the new row 2 with the "(in percent)" caption, one language per column
$ws.Range("A2").Value = "(пайыз менен)"
$ws.Range("B2").Value = "(в процентах)"
$ws.Range("C2").Value = "(in percent)"

# Style A2 like the other italic sub-caption cells (centered, vertical center)
$ws.Range("A2").Font.Name = "Times New Roman"
$ws.Range("A2").Font.Italic = $true
$ws.Range("A2").Font.Size = 9
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108

$ws.Range("B2:C2").HorizontalAlignment = -4108
$ws.Range("B2:C2").VerticalAlignment = -4108

# 3) Add the 2023 year column (T) of data
$ws.Range("T4").Value = 2023

$ws.Range("T5").Value = 16.899999999999999
$ws.Range("T6").Value = 1.6414476026646523
$ws.Range("T7").Value = 3.7332918174062506
$ws.Range("T8").Value = 1.4142735666882158
$ws.Range("T9").Value = 1.1221168654730374
$ws.Range("T10").Value = 3.7365451394949116
$ws.Range("T11").Value = 0.78756923781505217
$ws.Range("T12").Value = 2.0683169713107259
$ws.Range("T13").Value = 1.8858656419865651
$ws.Range("T14").Value = 0.49710203004505782

# Match number formatting of the neighbouring "S" column for the new T column cells
$ws.Range("T4").NumberFormat = $ws.Range("S4").NumberFormat
$ws.Range("T5:T14").NumberFormat = $ws.Range("S5").NumberFormat

# 4) Narrower first three columns
$ws.Columns.Item(1).ColumnWidth = 33.28515625
$ws.Columns.Item(2).ColumnWidth = 33.28515625
$ws.Columns.Item(3).ColumnWidth = 33.28515625

# 5) Taller header row to fit the extra caption line
$ws.Rows.Item(1).RowHeight = 45.75
